# Update statistics (想去人数 / 最低票价) for the latest data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 923
$ws.Range("F3").Value  = 552
$ws.Range("G3").Value  = 75
$ws.Range("F6").Value  = 710
$ws.Range("F9").Value  = 119
$ws.Range("F11").Value = 189
$ws.Range("F12").Value = 4876
$ws.Range("F13").Value = 35
$ws.Range("F15").Value = 470
$ws.Range("G15").Value = 58
$ws.Range("F17").Value = 525
$ws.Range("F18").Value = 317
$ws.Range("F22").Value = 694
$ws.Range("F24").Value = 291
$ws.Range("F25").Value = 988
$ws.Range("F27").Value = 1680
$ws.Range("F28").Value = 406

# --- Sheet: 演出 (Shows) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 70

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 159

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 923
$ws.Range("F6").Value  = 159
$ws.Range("F7").Value  = 552
$ws.Range("G7").Value  = 75
$ws.Range("F10").Value = 710
$ws.Range("F14").Value = 119
$ws.Range("F16").Value = 189
$ws.Range("F17").Value = 4876
$ws.Range("F18").Value = 35
$ws.Range("F21").Value = 470
$ws.Range("G21").Value = 58
$ws.Range("F23").Value = 525
$ws.Range("F24").Value = 317
$ws.Range("F31").Value = 70
$ws.Range("F32").Value = 694
$ws.Range("F37").Value = 291
$ws.Range("F38").Value = 988
$ws.Range("F40").Value = 1680
$ws.Range("F41").Value = 406
